$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data update (single-record financial refresh: report period changed
# from 2020-06-30 to 2019-12-31, values refreshed accordingly)

# J2 holds a text code ("001") that looks numeric, so force Text format
# while writing it, then clear the format back off so no stray number
# formatting is left behind on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 57989273.63
$ws.Range("P2").Value = 172467827.82
$ws.Range("Q2").Value = 122038275.78
$ws.Range("R2").Value = 18.2949211993
$ws.Range("S2").Value = 57790619.69
$ws.Range("T2").Value = 57790619.69
$ws.Range("U2").Value = 13.3342526808
$ws.Range("V2").Value = 29080985.14
$ws.Range("W2").Value = 5745406.17
$ws.Range("X2").Value = -287735.06
$ws.Range("Y2").Value = 62504171.59
$ws.Range("Z2").Value = 62508389.64
$ws.Range("AA2").Value = 4519116.01
$ws.Range("AG2").Value = 1837821.37
$ws.Range("AP2").Value = 17.3241604531
$ws.Range("AQ2").Value = 11.377811589009
$ws.Range("AR2").Value = 12.864376015857
$ws.Range("AS2").Value = 54166972.92
$ws.Range("AT2").Value = 11.304243790688
